$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb1"
$ws.Range("C2").Value = "Ephb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.079689000000001
$ws.Range("H2").Value = 21.239067
$ws.Range("I2").Value = 0.5033576067109902
$ws.Range("J2").Value = 0.5033576067109902
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.143611
$ws.Range("N2").Value = 0.430833
$ws.Range("O2").Value = 0.006856337892517759
$ws.Range("P2").Value = 0.006856337892517758
$ws.Range("Q2").Value = 1.016721216979
$ws.Range("R2").Value = 9.150490952811001
$ws.Range("S2").Value = 0.003451189832379613
$ws.Range("T2").Value = 0.003451189832379613

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb1"
$ws.Range("C3").Value = "Ephb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.079689000000001
$ws.Range("H3").Value = 21.239067
$ws.Range("I3").Value = 0.5033576067109902
$ws.Range("J3").Value = 0.5033576067109902
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.192962
$ws.Range("N3").Value = 51.578886
$ws.Range("O3").Value = 0.8208337581746376
$ws.Range("P3").Value = 0.8208337581746377
$ws.Range("Q3").Value = 121.720823948818
$ws.Range("R3").Value = 1095.487415539362
$ws.Range("S3").Value = 0.4131729160223732
$ws.Range("T3").Value = 0.4131729160223733

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb1"
$ws.Range("C4").Value = "Ephb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.079689000000001
$ws.Range("H4").Value = 21.239067
$ws.Range("I4").Value = 0.5033576067109902
$ws.Range("J4").Value = 0.5033576067109902
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.609156666666667
$ws.Range("N4").Value = 10.82747
$ws.Range("O4").Value = 0.1723099039328446
$ws.Range("P4").Value = 0.1723099039328446
$ws.Range("Q4").Value = 25.55170675227667
$ws.Range("R4").Value = 229.96536077049
$ws.Range("S4").Value = 0.08673350085623729
$ws.Range("T4").Value = 0.08673350085623729

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb1"
$ws.Range("C5").Value = "Ephb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.058683666666667
$ws.Range("H5").Value = 12.176051
$ws.Range("I5").Value = 0.2885676612136944
$ws.Range("J5").Value = 0.2885676612136945
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.143611
$ws.Range("N5").Value = 0.430833
$ws.Range("O5").Value = 0.006856337892517759
$ws.Range("P5").Value = 0.006856337892517758
$ws.Range("Q5").Value = 0.5828716200536668
$ws.Range("R5").Value = 5.245844580483
$ws.Range("S5").Value = 0.00197851739013468
$ws.Range("T5").Value = 0.001978517390134681

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb1"
$ws.Range("C6").Value = "Ephb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.058683666666667
$ws.Range("H6").Value = 12.176051
$ws.Range("I6").Value = 0.2885676612136944
$ws.Range("J6").Value = 0.2885676612136945
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.192962
$ws.Range("N6").Value = 51.578886
$ws.Range("O6").Value = 0.8208337581746376
$ws.Range("P6").Value = 0.8208337581746377
$ws.Range("Q6").Value = 69.78079405102066
$ws.Range("R6").Value = 628.027146459186
$ws.Range("S6").Value = 0.2368660778417024
$ws.Range("T6").Value = 0.2368660778417025

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb1"
$ws.Range("C7").Value = "Ephb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.058683666666667
$ws.Range("H7").Value = 12.176051
$ws.Range("I7").Value = 0.2885676612136944
$ws.Range("J7").Value = 0.2885676612136945
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.609156666666667
$ws.Range("N7").Value = 10.82747
$ws.Range("O7").Value = 0.1723099039328446
$ws.Range("P7").Value = 0.1723099039328446
$ws.Range("Q7").Value = 14.64842521344111
$ws.Range("R7").Value = 131.83582692097
$ws.Range("S7").Value = 0.04972306598185734
$ws.Range("T7").Value = 0.04972306598185735

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Efnb1"
$ws.Range("C8").Value = "Ephb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.463315
$ws.Range("H8").Value = 1.389945
$ws.Range("I8").Value = 0.03294115455541936
$ws.Range("J8").Value = 0.03294115455541936
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.143611
$ws.Range("N8").Value = 0.430833
$ws.Range("O8").Value = 0.006856337892517759
$ws.Range("P8").Value = 0.006856337892517758
$ws.Range("Q8").Value = 0.06653713046500001
$ws.Range("R8").Value = 0.5988341741850001
$ws.Range("S8").Value = 0.0002258556862016058
$ws.Range("T8").Value = 0.0002258556862016057

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Efnb1"
$ws.Range("C9").Value = "Ephb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.463315
$ws.Range("H9").Value = 1.389945
$ws.Range("I9").Value = 0.03294115455541936
$ws.Range("J9").Value = 0.03294115455541936
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.192962
$ws.Range("N9").Value = 51.578886
$ws.Range("O9").Value = 0.8208337581746376
$ws.Range("P9").Value = 0.8208337581746377
$ws.Range("Q9").Value = 7.965757189029999
$ws.Range("R9").Value = 71.69181470126999
$ws.Range("S9").Value = 0.02703921169233646
$ws.Range("T9").Value = 0.02703921169233646

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Efnb1"
$ws.Range("C10").Value = "Ephb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.463315
$ws.Range("H10").Value = 1.389945
$ws.Range("I10").Value = 0.03294115455541936
$ws.Range("J10").Value = 0.03294115455541936
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.609156666666667
$ws.Range("N10").Value = 10.82747
$ws.Range("O10").Value = 0.1723099039328446
$ws.Range("P10").Value = 0.1723099039328446
$ws.Range("Q10").Value = 1.672176421016667
$ws.Range("R10").Value = 15.04958778915
$ws.Range("S10").Value = 0.005676087176881298
$ws.Range("T10").Value = 0.005676087176881298

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Efnb1"
$ws.Range("C11").Value = "Ephb3"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.463241333333333
$ws.Range("H11").Value = 7.389724
$ws.Range("I11").Value = 0.175133577519896
$ws.Range("J11").Value = 0.175133577519896
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.143611
$ws.Range("N11").Value = 0.430833
$ws.Range("O11").Value = 0.006856337892517759
$ws.Range("P11").Value = 0.006856337892517758
$ws.Range("Q11").Value = 0.3537485511213334
$ws.Range("R11").Value = 3.183736960092
$ws.Range("S11").Value = 0.001200774983801859
$ws.Range("T11").Value = 0.001200774983801859

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Efnb1"
$ws.Range("C12").Value = "Ephb3"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.463241333333333
$ws.Range("H12").Value = 7.389724
$ws.Range("I12").Value = 0.175133577519896
$ws.Range("J12").Value = 0.175133577519896
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 17.192962
$ws.Range("N12").Value = 51.578886
$ws.Range("O12").Value = 0.8208337581746376
$ws.Range("P12").Value = 0.8208337581746377
$ws.Range("Q12").Value = 42.35041464082933
$ws.Range("R12").Value = 381.153731767464
$ws.Range("S12").Value = 0.1437555526182254
$ws.Range("T12").Value = 0.1437555526182255

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Efnb1"
$ws.Range("C13").Value = "Ephb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.463241333333333
$ws.Range("H13").Value = 7.389724
$ws.Range("I13").Value = 0.175133577519896
$ws.Range("J13").Value = 0.175133577519896
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.609156666666667
$ws.Range("N13").Value = 10.82747
$ws.Range("O13").Value = 0.1723099039328446
$ws.Range("P13").Value = 0.1723099039328446
$ws.Range("Q13").Value = 8.89022387980889
$ws.Range("R13").Value = 80.01201491828
$ws.Range("S13").Value = 0.03017724991786867
$ws.Range("T13").Value = 0.03017724991786867
